$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 9 (1-indexed): Mushroom / 70 000
$t.Cell(9, 1).Range.Text = "Mushroom"
$t.Cell(9, 3).Range.Text = "70 000"

# Row 10 (1-indexed): Leaf / 80 000 (trailing space)
$t.Cell(10, 1).Range.Text = "Leaf"
$t.Cell(10, 3).Range.Text = "80 000 "
